$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 98
$ws.Cells.Item(98, 8).Value = 744.5  # H98 was 696.1539
$ws.Cells.Item(98, 9).Value = 716.1111  # I98 was 670.8333
$ws.Cells.Item(98, 11).Value = 716.1111  # K98 was 670.8333
$ws.Cells.Item(98, 13).Value = 781.8889  # M98 was 827.1667
# Row 113
$ws.Cells.Item(113, 8).Value = 92991.37  # H113 was 51990.3
$ws.Cells.Item(113, 9).Value = 144986.42  # I113 was 127125.625
$ws.Cells.Item(113, 10).Value = 2000  # J113 was 1900.0834
$ws.Cells.Item(113, 11).Value = 144986.42  # K113 was 127125.625
$ws.Cells.Item(113, 12).Value = 2000  # L113 was 1900.0834
$ws.Cells.Item(113, 13).Value = -141732.42  # M113 was -123871.625
$ws.Cells.Item(113, 14).Value = -8508  # N113 was -8408.0834
# Row 122
$ws.Cells.Item(122, 8).Value = 744.5  # H122 was 696.1539
$ws.Cells.Item(122, 9).Value = 716.1111  # I122 was 670.8333
$ws.Cells.Item(122, 11).Value = 2148.3333  # K122 was 2012.4999
$ws.Cells.Item(122, 13).Value = 301.6667000000002  # M122 was 437.5001
# Row 124
$ws.Cells.Item(124, 8).Value = 35959.2  # H124 was 46551
$ws.Cells.Item(124, 10).Value = 35959.2  # J124 was 46551
$ws.Cells.Item(124, 12).Value = 35959.2  # L124 was 46551
$ws.Cells.Item(124, 14).Value = -45779.2  # N124 was -56371
# Row 129
$ws.Cells.Item(129, 8).Value = 946.84375  # H129 was 943.34326
$ws.Cells.Item(129, 10).Value = 966.2373  # J129 was 961.5161000000001
$ws.Cells.Item(129, 12).Value = 2898.7119  # L129 was 2884.5483
$ws.Cells.Item(129, 14).Value = -12898.7119  # N129 was -12884.5483
# Row 132
$ws.Cells.Item(132, 8).Value = 5560533.5  # H132 was 6103015
$ws.Cells.Item(132, 9).Value = 6950505  # I132 was 7149085
$ws.Cells.Item(132, 10).Value = 648.3333  # J132 was 940
$ws.Cells.Item(132, 11).Value = 20851515  # K132 was 21447255
$ws.Cells.Item(132, 12).Value = 1944.9999  # L132 was 2820
$ws.Cells.Item(132, 13).Value = -20848985  # M132 was -21444725
$ws.Cells.Item(132, 14).Value = -7004.9999  # N132 was -7880
# Row 137
$ws.Cells.Item(137, 8).Value = 1247.5424  # H137 was 1317.5178
$ws.Cells.Item(137, 9).Value = 954.85  # I137 was 1160.8125
$ws.Cells.Item(137, 10).Value = 1397.641  # J137 was 1380.2
$ws.Cells.Item(137, 11).Value = 2864.55  # K137 was 3482.4375
$ws.Cells.Item(137, 12).Value = 4192.923000000001  # L137 was 4140.6
$ws.Cells.Item(137, 13).Value = -314.5500000000002  # M137 was -932.4375
$ws.Cells.Item(137, 14).Value = -9292.923000000001  # N137 was -9240.6
# Row 138
$ws.Cells.Item(138, 8).Value = 4623.4653  # H138 was 4687.035
$ws.Cells.Item(138, 9).Value = 2771.7334  # I138 was 2898.2856
$ws.Cells.Item(138, 11).Value = 8315.200199999999  # K138 was 8694.856800000001
$ws.Cells.Item(138, 13).Value = -3175.200199999999  # M138 was -3554.856800000001

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 25527.225  # H32 was 26649.133
$ws.Cells.Item(32, 9).Value = 4145.9834  # I32 was 4355.3965
$ws.Cells.Item(32, 11).Value = 4145.9834  # K32 was 4355.3965
$ws.Cells.Item(32, 13).Value = -3858.9834  # M32 was -4068.3965
# Row 74
$ws.Cells.Item(74, 8).Value = 2236.3225  # H74 was 1588.3062
$ws.Cells.Item(74, 9).Value = 1398.1538  # I74 was 928.52
$ws.Cells.Item(74, 10).Value = 2841.6667  # J74 was 2275.5833
$ws.Cells.Item(74, 11).Value = 1398.1538  # K74 was 928.52
$ws.Cells.Item(74, 12).Value = 2841.6667  # L74 was 2275.5833
$ws.Cells.Item(74, 13).Value = -524.1538  # M74 was -54.51999999999998
$ws.Cells.Item(74, 14).Value = -4589.6667  # N74 was -4023.5833
# Row 77
$ws.Cells.Item(77, 8).Value = 2236.3225  # H77 was 1588.3062
$ws.Cells.Item(77, 9).Value = 1398.1538  # I77 was 928.52
$ws.Cells.Item(77, 10).Value = 2841.6667  # J77 was 2275.5833
$ws.Cells.Item(77, 11).Value = 6990.769  # K77 was 4642.6
$ws.Cells.Item(77, 12).Value = 14208.3335  # L77 was 11377.9165
$ws.Cells.Item(77, 13).Value = -2622.769  # M77 was -274.6000000000004
$ws.Cells.Item(77, 14).Value = -22944.3335  # N77 was -20113.9165
# Row 110
$ws.Cells.Item(110, 8).Value = 17277328  # H110 was 13917936
$ws.Cells.Item(110, 9).Value = 25050968  # I110 was 18556404
$ws.Cells.Item(110, 10).Value = 2573.889  # J110 was 2533.5557
$ws.Cells.Item(110, 11).Value = 25050968  # K110 was 18556404
$ws.Cells.Item(110, 12).Value = 2573.889  # L110 was 2533.5557
$ws.Cells.Item(110, 13).Value = -25048923  # M110 was -18554359
$ws.Cells.Item(110, 14).Value = -6663.889  # N110 was -6623.5557
# Row 122
$ws.Cells.Item(122, 8).Value = 3935.1333  # H122 was 2943.9546
$ws.Cells.Item(122, 9).Value = 4492.1816  # I122 was 3064.111
$ws.Cells.Item(122, 11).Value = 13476.5448  # K122 was 9192.332999999999
$ws.Cells.Item(122, 13).Value = -11026.5448  # M122 was -6742.332999999999
# Row 132
$ws.Cells.Item(132, 8).Value = 3440.8918  # H132 was 3112.442
$ws.Cells.Item(132, 9).Value = 5564.1177  # I132 was 4396.174
$ws.Cells.Item(132, 11).Value = 16692.3531  # K132 was 13188.522
$ws.Cells.Item(132, 13).Value = -14162.3531  # M132 was -10658.522

$ws = $wb.Worksheets.Item("BSM")
# Row 8
$ws.Cells.Item(8, 8).Value = 1100  # H8 was 27416.666
$ws.Cells.Item(8, 9).Value = 1100  # I8 was 950
$ws.Cells.Item(8, 10).Value = 0  # J8 was 40650
$ws.Cells.Item(8, 11).Value = 1100  # K8 was 950
$ws.Cells.Item(8, 12).Value = 0  # L8 was 40650
$ws.Cells.Item(8, 13).Value = -960  # M8 was -810
$ws.Cells.Item(8, 14).ClearContents()  # N8 was -40930
# Row 134
$ws.Cells.Item(134, 8).Value = 2689.08  # H134 was 1765.674
$ws.Cells.Item(134, 9).Value = 3466.4666  # I134 was 1956.7812
$ws.Cells.Item(134, 10).Value = 1523  # J134 was 1328.8572
$ws.Cells.Item(134, 11).Value = 10399.3998  # K134 was 5870.3436
$ws.Cells.Item(134, 12).Value = 4569  # L134 was 3986.5716
$ws.Cells.Item(134, 13).Value = -7864.399800000001  # M134 was -3335.3436
$ws.Cells.Item(134, 14).Value = -9639  # N134 was -9056.571599999999

$ws = $wb.Worksheets.Item("CRP")
# Row 99
$ws.Cells.Item(99, 8).Value = 9188  # H99 was 9599.235000000001
$ws.Cells.Item(99, 10).Value = 14322.777  # J99 was 15838.5
$ws.Cells.Item(99, 12).Value = 14322.777  # L99 was 15838.5
$ws.Cells.Item(99, 14).Value = -17318.777  # N99 was -18834.5
# Row 126
$ws.Cells.Item(126, 8).Value = 9188  # H126 was 9599.235000000001
$ws.Cells.Item(126, 10).Value = 14322.777  # J126 was 15838.5
$ws.Cells.Item(126, 12).Value = 42968.331  # L126 was 47515.5
$ws.Cells.Item(126, 14).Value = -47908.331  # N126 was -52455.5
# Row 132
$ws.Cells.Item(132, 8).Value = 2355.606  # H132 was 2809.2222
$ws.Cells.Item(132, 9).Value = 2087.4827  # I132 was 2676.238
$ws.Cells.Item(132, 10).Value = 4299.5  # J132 was 3274.6667
$ws.Cells.Item(132, 11).Value = 6262.4481  # K132 was 8028.714
$ws.Cells.Item(132, 12).Value = 12898.5  # L132 was 9824.000100000001
$ws.Cells.Item(132, 13).Value = -3732.4481  # M132 was -5498.714
$ws.Cells.Item(132, 14).Value = -17958.5  # N132 was -14884.0001
# Row 134
$ws.Cells.Item(134, 8).Value = 2098.4285  # H134 was 1487.75
$ws.Cells.Item(134, 9).Value = 1578.5454  # I134 was 1222.75
$ws.Cells.Item(134, 10).Value = 4004.6667  # J134 was 2150.25
$ws.Cells.Item(134, 11).Value = 4735.6362  # K134 was 3668.25
$ws.Cells.Item(134, 12).Value = 12014.0001  # L134 was 6450.75
$ws.Cells.Item(134, 13).Value = -2200.6362  # M134 was -1133.25
$ws.Cells.Item(134, 14).Value = -17084.0001  # N134 was -11520.75

$ws = $wb.Worksheets.Item("CUL")
# Row 44
$ws.Cells.Item(44, 8).Value = 850.5  # H44 was 900.6
$ws.Cells.Item(44, 9).Value = 467.33334  # I44 was 499.5
$ws.Cells.Item(44, 10).Value = 2000  # J44 was 1168
$ws.Cells.Item(44, 11).Value = 1402.00002  # K44 was 1498.5
$ws.Cells.Item(44, 12).Value = 6000  # L44 was 3504
$ws.Cells.Item(44, 13).Value = -1004.00002  # M44 was -1100.5
$ws.Cells.Item(44, 14).Value = -6796  # N44 was -4300
# Row 47
$ws.Cells.Item(47, 8).Value = 1714.3334  # H47 was 0
$ws.Cells.Item(47, 9).Value = 143  # I47 was 0
$ws.Cells.Item(47, 10).Value = 2500  # J47 was 0
$ws.Cells.Item(47, 11).Value = 429  # K47 was 0
$ws.Cells.Item(47, 12).Value = 7500  # L47 was 0
$ws.Cells.Item(47, 13).Value = 2  # M47 was None
$ws.Cells.Item(47, 14).Value = -8362  # N47 was None
# Row 129
$ws.Cells.Item(129, 8).Value = 1572.7059  # H129 was 1875.5333
$ws.Cells.Item(129, 9).Value = 504.2857  # I129 was 533.3333
$ws.Cells.Item(129, 10).Value = 2320.6  # J129 was 2770.3333
$ws.Cells.Item(129, 11).Value = 1512.8571  # K129 was 1599.9999
$ws.Cells.Item(129, 12).Value = 6961.799999999999  # L129 was 8310.999899999999
$ws.Cells.Item(129, 13).Value = 3487.1429  # M129 was 3400.0001
$ws.Cells.Item(129, 14).Value = -16961.8  # N129 was -18310.9999
# Row 131
$ws.Cells.Item(131, 8).Value = 1201.24  # H131 was 1245.806
$ws.Cells.Item(131, 10).Value = 1253.7384  # J131 was 1313.4912
$ws.Cells.Item(131, 12).Value = 3761.2152  # L131 was 3940.4736
$ws.Cells.Item(131, 14).Value = -13841.2152  # N131 was -14020.4736

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Cells.Item(102, 8).Value = 195432.27  # H102 was 232626.92
$ws.Cells.Item(102, 9).Value = 1987.08  # I102 was 1978.85
$ws.Cells.Item(102, 11).Value = 1987.08  # K102 was 1978.85
$ws.Cells.Item(102, 13).Value = -365.0799999999999  # M102 was -356.8499999999999
# Row 113
$ws.Cells.Item(113, 8).Value = 1640.3846  # H113 was 1692.6666
$ws.Cells.Item(113, 10).Value = 1875.5  # J113 was 1998.7142
$ws.Cells.Item(113, 12).Value = 1875.5  # L113 was 1998.7142
$ws.Cells.Item(113, 14).Value = -6215.5  # N113 was -6338.7142
# Row 122
$ws.Cells.Item(122, 8).Value = 6493.4614  # H122 was 5337.4736
$ws.Cells.Item(122, 9).Value = 6343.857  # I122 was 5040.5
$ws.Cells.Item(122, 10).Value = 6668  # J122 was 5667.4443
$ws.Cells.Item(122, 11).Value = 19031.571  # K122 was 15121.5
$ws.Cells.Item(122, 12).Value = 20004  # L122 was 17002.3329
$ws.Cells.Item(122, 13).Value = -16581.571  # M122 was -12671.5
$ws.Cells.Item(122, 14).Value = -24904  # N122 was -21902.3329
# Row 132
$ws.Cells.Item(132, 8).Value = 2539.75  # H132 was 2529.5334
$ws.Cells.Item(132, 9).Value = 2149.0557  # I132 was 2080.75
$ws.Cells.Item(132, 10).Value = 3042.0715  # J132 was 3042.4285
$ws.Cells.Item(132, 11).Value = 6447.1671  # K132 was 6242.25
$ws.Cells.Item(132, 12).Value = 9126.2145  # L132 was 9127.2855
$ws.Cells.Item(132, 13).Value = -3917.1671  # M132 was -3712.25
$ws.Cells.Item(132, 14).Value = -14186.2145  # N132 was -14187.2855

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Cells.Item(40, 8).Value = 85764.336  # H40 was 79205.84
$ws.Cells.Item(40, 9).Value = 334733.34  # I40 was 201439.8
$ws.Cells.Item(40, 10).Value = 2774.6667  # J40 was 2809.625
$ws.Cells.Item(40, 11).Value = 334733.34  # K40 was 201439.8
$ws.Cells.Item(40, 12).Value = 2774.6667  # L40 was 2809.625
$ws.Cells.Item(40, 13).Value = -334597.34  # M40 was -201303.8
$ws.Cells.Item(40, 14).Value = -3046.6667  # N40 was -3081.625
# Row 132
$ws.Cells.Item(132, 8).Value = 4006.077  # H132 was 4728.7617
$ws.Cells.Item(132, 9).Value = 5662.25  # I132 was 6649.3
$ws.Cells.Item(132, 10).Value = 2586.5  # J132 was 2982.818
$ws.Cells.Item(132, 11).Value = 16986.75  # K132 was 19947.9
$ws.Cells.Item(132, 12).Value = 7759.5  # L132 was 8948.454000000002
$ws.Cells.Item(132, 13).Value = -14456.75  # M132 was -17417.9
$ws.Cells.Item(132, 14).Value = -12819.5  # N132 was -14008.454

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Cells.Item(122, 8).Value = 1242.8572  # H122 was 911.85187
$ws.Cells.Item(122, 9).Value = 1242.8572  # I122 was 936.5294
$ws.Cells.Item(122, 10).Value = 0  # J122 was 869.9
$ws.Cells.Item(122, 11).Value = 3728.5716  # K122 was 2809.5882
$ws.Cells.Item(122, 12).Value = 0  # L122 was 2609.7
$ws.Cells.Item(122, 13).Value = -1278.5716  # M122 was -359.5882000000001
$ws.Cells.Item(122, 14).ClearContents()  # N122 was -7509.7
# Row 123
$ws.Cells.Item(123, 8).Value = 35517  # H123 was 37330
$ws.Cells.Item(123, 10).Value = 35517  # J123 was 37330
$ws.Cells.Item(123, 12).Value = 35517  # L123 was 37330
$ws.Cells.Item(123, 14).Value = -45317  # N123 was -47130
